$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("B2").Value = "<i>"
$ws.Range("C2").Value = 40

# Row 3
$ws.Range("C3").Value = 45

# Row 4
$ws.Range("B4").Value = "<he>"
$ws.Range("C4").Value = 42

# Row 5
$ws.Range("C5").Value = 38

# Row 6
$ws.Range("C6").Value = 31

# Row 7
$ws.Range("C7").Value = 42

# Row 8
$ws.Range("B8").Value = "<novem>"
$ws.Range("C8").Value = 39

# Row 9
$ws.Range("C9").Value = 42

# Row 10
$ws.Range("B10").Value = "<tab>"
$ws.Range("C10").Value = 39

# Row 11
$ws.Range("B11").Value = "<eight>"

# Row 12
$ws.Range("B12").Value = "<and>"

# Row 13
$ws.Range("B13").Value = "<seven>"
$ws.Range("C13").Value = 42

# Row 15
$ws.Range("C15").Value = 40

# Row 16
$ws.Range("B16").Value = "<my>"
$ws.Range("C16").Value = 43

# Row 17
$ws.Range("B17").Value = "<enten>"
$ws.Range("C17").Value = 41

# Row 18
$ws.Range("B18").Value = "<which>"
$ws.Range("C18").Value = 31
